$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.424199999999995
$ws.Range("D10").Value = -8.162899999999997
$ws.Range("D12").Value = -7.988599999999999
$ws.Range("D18").Value = -8.071999999999992
$ws.Range("D25").Value = -8.268399999999996
